$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.021590333333334
$ws.Range("H2").Value = 6.064771
$ws.Range("I2").Value = 0.01116262347650641
$ws.Range("J2").Value = 0.01116262347650641
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.09934133333334
$ws.Range("N2").Value = 63.29802400000001
$ws.Range("O2").Value = 0.2917236204149438
$ws.Range("P2").Value = 0.2917236204149438
$ws.Range("Q2").Value = 42.65422447916712
$ws.Range("R2").Value = 383.8880203125041
$ws.Range("S2").Value = 0.003256400933895297
$ws.Range("T2").Value = 0.003256400933895297
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.021590333333334
$ws.Range("H3").Value = 6.064771
$ws.Range("I3").Value = 0.01116262347650641
$ws.Range("J3").Value = 0.01116262347650641
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.81943766666667
$ws.Range("N3").Value = 107.458313
$ws.Range("O3").Value = 0.4952465516465762
$ws.Range("P3").Value = 0.4952465516465762
$ws.Range("Q3").Value = 72.41222893236923
$ws.Range("R3").Value = 651.7100603913231
$ws.Range("S3").Value = 0.005528250784068917
$ws.Range("T3").Value = 0.005528250784068917
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.021590333333334
$ws.Range("H4").Value = 6.064771
$ws.Range("I4").Value = 0.01116262347650641
$ws.Range("J4").Value = 0.01116262347650641
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 15.40769666666667
$ws.Range("N4").Value = 46.22309
$ws.Range("O4").Value = 0.2130298279384801
$ws.Range("P4").Value = 0.2130298279384801
$ws.Range("Q4").Value = 31.14805064026556
$ws.Range("R4").Value = 280.33245576239
$ws.Range("S4").Value = 0.0023779717585422
$ws.Range("T4").Value = 0.0023779717585422
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 25.140634
$ws.Range("H5").Value = 75.421902
$ws.Range("I5").Value = 0.1388191398995883
$ws.Range("J5").Value = 0.1388191398995883
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.09934133333334
$ws.Range("N5").Value = 63.29802400000001
$ws.Range("O5").Value = 0.2917236204149438
$ws.Range("P5").Value = 0.2917236204149438
$ws.Range("Q5").Value = 530.4508181024054
$ws.Range("R5").Value = 4774.057362921649
$ws.Range("S5").Value = 0.04049682207439647
$ws.Range("T5").Value = 0.04049682207439647
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 25.140634
$ws.Range("H6").Value = 75.421902
$ws.Range("I6").Value = 0.1388191398995883
$ws.Range("J6").Value = 0.1388191398995883
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.81943766666667
$ws.Range("N6").Value = 107.458313
$ws.Range("O6").Value = 0.4952465516465762
$ws.Range("P6").Value = 0.4952465516465762
$ws.Range("Q6").Value = 900.5233724634808
$ws.Range("R6").Value = 8104.710352171327
$ws.Range("S6").Value = 0.06874970033781475
$ws.Range("T6").Value = 0.06874970033781473
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 25.140634
$ws.Range("H7").Value = 75.421902
$ws.Range("I7").Value = 0.1388191398995883
$ws.Range("J7").Value = 0.1388191398995883
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.40769666666667
$ws.Range("N7").Value = 46.22309
$ws.Range("O7").Value = 0.2130298279384801
$ws.Range("P7").Value = 0.2130298279384801
$ws.Range("Q7").Value = 387.3592626796867
$ws.Range("R7").Value = 3486.23336411718
$ws.Range("S7").Value = 0.02957261748737709
$ws.Range("T7").Value = 0.02957261748737709
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 153.9412893333333
$ws.Range("H8").Value = 461.823868
$ws.Range("I8").Value = 0.8500182366239053
$ws.Range("J8").Value = 0.8500182366239052
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.09934133333334
$ws.Range("N8").Value = 63.29802400000001
$ws.Range("O8").Value = 0.2917236204149438
$ws.Range("P8").Value = 0.2917236204149438
$ws.Range("Q8").Value = 3248.059808937426
$ws.Range("R8").Value = 29232.53828043683
$ws.Range("S8").Value = 0.247970397406652
$ws.Range("T8").Value = 0.247970397406652
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 153.9412893333333
$ws.Range("H9").Value = 461.823868
$ws.Range("I9").Value = 0.8500182366239053
$ws.Range("J9").Value = 0.8500182366239052
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 35.81943766666667
$ws.Range("N9").Value = 107.458313
$ws.Range("O9").Value = 0.4952465516465762
$ws.Range("P9").Value = 0.4952465516465762
$ws.Range("Q9").Value = 5514.090417601632
$ws.Range("R9").Value = 49626.81375841469
$ws.Range("S9").Value = 0.4209686005246925
$ws.Range("T9").Value = 0.4209686005246925
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 153.9412893333333
$ws.Range("H10").Value = 461.823868
$ws.Range("I10").Value = 0.8500182366239053
$ws.Range("J10").Value = 0.8500182366239052
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.40769666666667
$ws.Range("N10").Value = 46.22309
$ws.Range("O10").Value = 0.2130298279384801
$ws.Range("P10").Value = 0.2130298279384801
$ws.Range("Q10").Value = 2371.880690523569
$ws.Range("R10").Value = 21346.92621471212
$ws.Range("S10").Value = 0.1810792386925608
$ws.Range("T10").Value = 0.1810792386925608
